$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 24: Excel auto-fit removed the custom 30pt row height ---
$ws.Rows.Item(24).AutoFit()

# --- New row 71: "Création de vidéos explicatives " / 1h ---
$ws.Range("B70:E70").Copy($ws.Range("B71"))
$ws.Range("B71").Value = 52
$ws.Range("C71").Value = 44967
$ws.Range("D71").Value = "Création de vidéos explicatives "
$ws.Range("E71").Value = "1h"

# --- New row 72: meeting note with Leonel / 2h30 ---
$ws.Range("B71:C71").Copy($ws.Range("B72"))
$ws.Range("B72").Value = 53
$ws.Range("C72").Value = 44967
$ws.Range("D72").Value = "réunion avec Leonel pour finaliser les détails de l'application et sa présentation`n"
$ws.Range("D72").HorizontalAlignment = -4131
$ws.Range("D72").VerticalAlignment = -4160
$ws.Range("E72").Value = "2h30"
$ws.Rows.Item(72).RowHeight = 16.5

# --- Update the view: scroll/selection moved to G25 ---
[void]$ws.Range("G25").Select()
